$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C1 with the new message text (adds a new shared string entry)
$ws.Range("C1").Value = "Epic sadface: Username and password do not match any user in this service1"

# Update the selection to C1
$ws.Range("C1").Select()
